$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") and column E ("Volume(1h)") values refreshed to the latest
# scrape snapshot. Numeric-looking prices are entered with a leading apostrophe
# (exactly like typing '0.9978 into Excel) so they stay literal text instead of
# being reinterpreted/renormalised as numbers (e.g. "97.40" -> 97.4).

$ws.Range("D2").Value = "24.910.02"
$ws.Range("E2").Value = "  -3.82%  "
$ws.Range("D3").Value = "1.636.80"
$ws.Range("E3").Value = "  -5.99%  "
$ws.Range("D4").Value = "'0.9978"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "'235.65"
$ws.Range("E5").Value = "  -4.46%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.4710"
$ws.Range("E7").Value = "  -6.78%  "
$ws.Range("D8").Value = "'0.2551"
$ws.Range("E8").Value = "  -6.07%  "
$ws.Range("D9").Value = "'0.06008"
$ws.Range("E9").Value = "  -2.73%  "
$ws.Range("D10").Value = "'0.07136"
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("D11").Value = "1.635.89"
$ws.Range("E11").Value = "  -6.03%  "
$ws.Range("D12").Value = "'14.78"
$ws.Range("E12").Value = "  -1.96%  "
$ws.Range("D13").Value = "'0.6141"
$ws.Range("E13").Value = "  -4.92%  "
$ws.Range("D14").Value = "'4.412"
$ws.Range("E14").Value = "  -4.58%  "
$ws.Range("D15").Value = "'72.56"
$ws.Range("E15").Value = "  -6.33%  "
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "'0.9977"
$ws.Range("D18").Value = "24.904.15"
$ws.Range("E18").Value = "  -3.88%  "
$ws.Range("D19").Value = "'0.000006564"
$ws.Range("E19").Value = "  -3.55%  "
$ws.Range("E20").Value = "  -5.14%  "
$ws.Range("D21").Value = "'4.404"
$ws.Range("E21").Value = "  +2.98%  "
$ws.Range("D22").Value = "1.842.50"
$ws.Range("E22").Value = "  -6.27%  "
$ws.Range("D23").Value = "'8.554"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").Value = "'5.245"
$ws.Range("E24").Value = "  -2.26%  "
$ws.Range("D25").Value = "'132.70"
$ws.Range("E25").Value = "  -2.49%  "
$ws.Range("D26").Value = "'14.76"
$ws.Range("E26").Value = "  -3.08%  "
$ws.Range("D27").Value = "'1.371"
$ws.Range("E27").Value = "  -8.67%  "
$ws.Range("D28").Value = "'102.36"
$ws.Range("E28").Value = "  -2.76%  "
$ws.Range("D29").Value = "'1.651"
$ws.Range("E29").Value = "  -6.29%  "
$ws.Range("D30").Value = "'3.721"
$ws.Range("E30").Value = "  -4.82%  "
$ws.Range("D31").Value = "'0.07736"
$ws.Range("E31").Value = "  -5.91%  "
$ws.Range("E32").Value = "  -2.41%  "
$ws.Range("D33").Value = "'0.04367"
$ws.Range("E33").Value = "  -6.45%  "
$ws.Range("D34").Value = "'0.9989"
$ws.Range("D35").Value = "'2.598"
$ws.Range("E35").Value = "  -2.06%  "
$ws.Range("D36").Value = "'0.9175"
$ws.Range("E36").Value = "  -7.49%  "
$ws.Range("D37").Value = "'0.5798"
$ws.Range("E37").Value = "  -6.41%  "
$ws.Range("D38").Value = "'2.547"
$ws.Range("E38").Value = "  -6.71%  "
$ws.Range("D39").Value = "'0.01552"
$ws.Range("E39").Value = "  -2.88%  "
$ws.Range("D40").Value = "'0.9978"
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("D41").Value = "'0.8147"
$ws.Range("E41").Value = "  +7.45%  "
$ws.Range("D42").Value = "'1.794"
$ws.Range("E42").Value = "  -6.09%  "
$ws.Range("D43").Value = "'97.40"
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("D44").Value = "'0.3696"
$ws.Range("E44").Value = "  -3.99%  "
$ws.Range("D45").Value = "'4.728"
$ws.Range("E45").Value = "  -5.02%  "
$ws.Range("D46").Value = "'0.1126"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").Value = "'0.05211"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").Value = "'6.067"
$ws.Range("D49").Value = "'29.43"
$ws.Range("D50").Value = "'0.9997"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").Value = "'0.9998"
$ws.Range("E51").Value = "  -0.34%  "
